$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "在人物Cpp中实现" notes next to the first three top rows.
$ws.Range("D2").Value = "在人物Cpp中实现"
$ws.Range("D3").Value = "在人物Cpp中实现"
$ws.Range("D4").Value = "在人物Cpp中实现"

# Row 5 gets a new task entry plus a matching note.
$ws.Range("A5").Value = "已从第三人称模板中迁移动画"
$ws.Range("D5").Value = "在人物蓝图中实现"

# Insert a new row below row 6 (inside the merged B6:B17 block) so the
# merge grows to B6:B18 and the old "死亡动画" entry slides down to row 7.
$ws.Rows(7).Insert()

# Carry the previous row-6 detail text ("死亡动画") down into the freshly
# inserted row...
$ws.Range("C7").Value = "死亡动画"
# ...then replace row 6's detail text with the new task.
$ws.Range("C6").Value = "动画蓝图类的cpp创建"

# Match the author's final selection.
[void]$ws.Range("C6").Select()
